# Updates crypto price/volume figures per the scraper's latest run.
# Note: some Price (column D) values are plain numeric-looking strings
# (e.g. "96.87"). The source cells are plain text, so a leading apostrophe
# is used to force Excel to keep them as text instead of auto-converting
# them to numbers (this mirrors typing '96.87 directly into Excel).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.713.80'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '2.296.46'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''96.87'
$ws.Range('E5').Value = '  +2.43%  '
$ws.Range('D6').Value = '''268.60'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('D7').Value = '''0.624'
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.609'
$ws.Range('E9').Value = '  -2.20%  '
$ws.Range('D10').Value = '''45.34'
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').Value = '''7.90'
$ws.Range('E12').Value = '  -2.10%  '
$ws.Range('E13').Value = '  +1.53%  '
$ws.Range('D14').Value = '''15.65'
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('D15').Value = '2.640.63'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').Value = '''0.857'
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').Value = '2.296.87'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').Value = '43.743.24'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('E19').Value = '  +2.81%  '
$ws.Range('E20').Value = '  -2.45%  '
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('D22').Value = '''2.52'
$ws.Range('E22').Value = '  +10.18%  '
$ws.Range('D23').Value = '''233.15'
$ws.Range('E23').Value = '  -1.89%  '
$ws.Range('D24').Value = '''9.08'
$ws.Range('E24').Value = '  -5.19%  '
$ws.Range('E25').Value = '  +5.64%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = '''11.29'
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('E28').Value = '  +2.27%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.29'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '''38.88'
$ws.Range('E30').Value = '  +0.42%  '
$ws.Range('D31').Value = '''175.22'
$ws.Range('E31').Value = '  +1.96%  '
$ws.Range('D32').Value = '''21.87'
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('D33').Value = '''0.0904'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('D34').Value = '''5.42'
$ws.Range('E34').Value = '  -1.42%  '
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('D36').Value = '''4.58'
$ws.Range('E36').Value = '  +3.37%  '
$ws.Range('D37').Value = '''0.107'
$ws.Range('E37').Value = '  -1.28%  '
$ws.Range('D38').Value = '''0.0351'
$ws.Range('E38').Value = '  -1.17%  '
$ws.Range('E39').Value = '  -1.86%  '
$ws.Range('D40').Value = '''0.239'
$ws.Range('E40').Value = '  +2.53%  '
$ws.Range('D41').Value = '''2.32'
$ws.Range('E41').Value = '  +1.05%  '
$ws.Range('E42').Value = '  -1.17%  '
$ws.Range('D43').Value = '''12.20'
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('D44').Value = '''64.26'
$ws.Range('E44').Value = '  +3.87%  '
$ws.Range('D45').Value = '''8.81'
$ws.Range('E45').Value = '  -2.86%  '
$ws.Range('E46').Value = '  -5.05%  '
$ws.Range('E47').Value = '  -0.51%  '
$ws.Range('D48').Value = '''97.34'
$ws.Range('E48').Value = '  -3.16%  '
$ws.Range('E49').Value = '  -1.29%  '
$ws.Range('D50').Value = '''0.437'
$ws.Range('E50').Value = '  +3.76%  '
$ws.Range('E51').Value = '  +11.72%  '
